$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so Excel keeps them as text
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "26.929.92"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "1.864.03"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "304.94"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.5065"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").Value = "0.3649"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").Value = "0.07179"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "0.8966"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "20.85"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07485"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.845.14"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "92.66"
$ws.Range("D15").Value = "5.240"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "0.000008503"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "0.9989"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "26.972.51"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "5.041"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "2.086.95"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "10.39"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "6.405"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "147.41"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").Value = "1.790"
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("D27").Value = "17.89"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "2.081"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "113.19"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "4.698"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "0.09251"
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("D33").Value = "0.05111"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").Value = "0.7529"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").Value = "2.979"
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("D36").Value = "1.153"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "3.270"
$ws.Range("E37").Value = "  +7.14%  "
$ws.Range("D38").Value = "2.544"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("D40").Value = "0.5545"
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("D41").Value = "1.071"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").Value = "118.26"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").Value = "6.516"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("D44").Value = "8.530"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("D45").Value = "0.1473"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "0.4693"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "0.9985"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "10.09"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "1.566"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "36.84"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").Value = "62.95"
$ws.Range("E51").Value = "  -2.25%  "
